# Fruta / hortaliza, semanal
# The weekly refresh reorders the daily price rows (columns D, L, M, N, O, P, S
# carry the per-record date/quality/volume/price data; A,B,C,E-K,Q,R,T are
# constant for every row in this sheet) without changing the set of records.
# We snapshot the 7 varying columns for each data row, then re-write them back
# into their new row positions per the mapping below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maps the destination (new) row number -> source (old) row number.
$rowMap = @{
    2 = 13;  3 = 4;   4 = 5;   5 = 6;   6 = 2;   7 = 3;
    8 = 20;  9 = 21;  10 = 22; 11 = 27; 12 = 28; 13 = 29;
    14 = 16; 15 = 17; 16 = 18; 17 = 32; 18 = 33; 19 = 14;
    20 = 15; 21 = 7;  22 = 8;  23 = 19; 24 = 25; 25 = 26;
    26 = 9;  27 = 10; 28 = 23; 29 = 24; 30 = 11; 31 = 12;
    32 = 30; 33 = 31
}

$cols = @("D", "L", "M", "N", "O", "P", "S")

# Snapshot the current (pre-shuffle) values of the varying columns for every
# data row before any writes happen, so row N being overwritten doesn't
# clobber data still needed as the source for another destination row.
$snapshot = @{}
for ($r = 2; $r -le 33; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value()
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcVals[$col]
    }
}
